$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("kinetics_tyms")

# 1) Remove the empty "Sheet1" worksheet.
$wb.Worksheets.Item("Sheet1").Delete()

# 2) Re-point the table/body font to Arial for the data rows (A2:E5),
#    matching the kcat/Km table's new look.
$ws.Range("A2:E5").Font.Name = "Arial"

# 3) Recode the WT R166Q kcat/Km entries as asterisked placeholder text
#    (values too low to fit reliably) and right-align them.
$ws.Range("B5").Value = "0.001*"
$ws.Range("C5").Value = "10*"
$ws.Range("B5:C5").HorizontalAlignment = -4152   # xlRight

# 4) Swap the table's visual style.
$tbl = $ws.ListObjects.Item(1)
$tbl.TableStyle = "TableStyleLight1"

# 5) Update the selection to match the authored view.
$ws.Range("A2:E5").Select()
$ws.Range("C5").Activate()
